$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric (e.g. "216.12") but must remain
# stored as text, matching the source sheet where every Price/Volume
# cell is inline text. Mark them as Text format before assigning so
# Excel does not silently coerce them into numbers.
$textForceRefs = @("D5", "D9", "D10", "D11", "D14", "D16", "D18", "D19", "D23", "D24", "D25", "D26", "D27", "D28", "D31", "D37", "D38", "D40", "D43", "D44", "D46", "D47", "D50", "D51")
foreach ($ref in $textForceRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

# --- Price/Volume/Coin/Link updates ---

# Row 2
$ws.Range("D2").Value = "26.870.02"
$ws.Range("E2").Value = "  +0.19%  "

# Row 3
$ws.Range("D3").Value = "1.671.10"
$ws.Range("E3").Value = "  +2.03%  "

# Row 4
$ws.Range("E4").Value = "  +0.26%  "

# Row 5
$ws.Range("D5").Value = "216.12"
$ws.Range("E5").Value = "  +0.79%  "

# Row 6
$ws.Range("E6").Value = "  +5.63%  "

# Row 7
$ws.Range("E7").Value = "  +0.24%  "

# Row 8
$ws.Range("E8").Value = "  +3.05%  "

# Row 9
$ws.Range("D9").Value = "0.0620"
$ws.Range("E9").Value = "  +1.37%  "

# Row 10
$ws.Range("D10").Value = "20.33"
$ws.Range("E10").Value = "  +4.59%  "

# Row 11
$ws.Range("D11").Value = "0.0893"
$ws.Range("E11").Value = "  +3.85%  "

# Row 12
$ws.Range("D12").Value = "1.910.37"
$ws.Range("E12").Value = "  +2.29%  "

# Row 13
$ws.Range("D13").Value = "1.695.67"
$ws.Range("E13").Value = "  +4.06%  "

# Row 14
$ws.Range("D14").Value = "4.09"
$ws.Range("E14").Value = "  +0.93%  "

# Row 15
$ws.Range("E15").Value = "  +1.86%  "

# Row 16
$ws.Range("D16").Value = "65.69"
$ws.Range("E16").Value = "  +1.81%  "

# Row 17
$ws.Range("D17").Value = "26.895.15"
$ws.Range("E17").Value = "  +0.35%  "

# Row 18
$ws.Range("D18").Value = "232.75"
$ws.Range("E18").Value = "  -3.74%  "

# Row 19
$ws.Range("D19").Value = "7.86"
$ws.Range("E19").Value = "  -0.04%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0735"
$ws.Range("E20").Value = "  +1.09%  "

# Row 21
$ws.Range("E21").Value = "  +0.22%  "

# Row 22
$ws.Range("E22").Value = "  +2.08%  "

# Row 23
$ws.Range("D23").Value = "2.21"
$ws.Range("E23").Value = "  -1.30%  "

# Row 24
$ws.Range("D24").Value = "9.21"
$ws.Range("E24").Value = "  -0.19%  "

# Row 25
$ws.Range("D25").Value = "145.82"
$ws.Range("E25").Value = "  -0.09%  "

# Row 26
$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").Value = "0.116"
$ws.Range("E26").Value = "  +2.53%  "

# Row 27
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "7.15"
$ws.Range("E27").Value = "  +1.08%  "

# Row 28
$ws.Range("D28").Value = "15.98"
$ws.Range("E28").Value = "  +0.96%  "

# Row 29
$ws.Range("E29").Value = "  +0.20%  "

# Row 30
$ws.Range("E30").Value = "  +0.39%  "

# Row 31
$ws.Range("D31").Value = "1.17"
$ws.Range("E31").Value = "  +0.99%  "

# Row 32
$ws.Range("E32").Value = "  +1.98%  "

# Row 33
$ws.Range("D33").Value = "1.465.79"
$ws.Range("E33").Value = "  -2.98%  "

# Row 34
$ws.Range("E34").Value = "  +5.06%  "

# Row 35
$ws.Range("E35").Value = "  +3.64%  "

# Row 36
$ws.Range("E36").Value = "  +0.22%  "

# Row 37
$ws.Range("D37").Value = "0.904"
$ws.Range("E37").Value = "  +5.44%  "

# Row 38
$ws.Range("D38").Value = "0.570"
$ws.Range("E38").Value = "  -0.56%  "

# Row 39
$ws.Range("E39").Value = "  +0.56%  "

# Row 40
$ws.Range("D40").Value = "5.98"
$ws.Range("E40").Value = "  +0.44%  "

# Row 41
$ws.Range("E41").Value = "  +0.28%  "

# Row 42
$ws.Range("E42").Value = "  +4.40%  "

# Row 43
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "0.973"
$ws.Range("E43").Value = "  +6.77%  "

# Row 44
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "65.84"
$ws.Range("E44").Value = "  +2.58%  "

# Row 45
$ws.Range("D45").Value = "1.818.32"
$ws.Range("E45").Value = "  +2.27%  "

# Row 46
$ws.Range("D46").Value = "0.780"
$ws.Range("E46").Value = "  +2.02%  "

# Row 47
$ws.Range("D47").Value = "90.45"
$ws.Range("E47").Value = "  +0.12%  "

# Row 48
$ws.Range("E48").Value = "  +0.65%  "

# Row 49
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0103"
$ws.Range("E49").Value = "  +16.74%  "

# Row 50
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.1000"
$ws.Range("E50").Value = "  +2.39%  "

# Row 51
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.0508"
$ws.Range("E51").Value = "  +1.38%  "
